# Generate Report for Handoff
# Adds a new localization-status row (b16ec6b2-f807-46c9-9097-3be3590d0583)
# to the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$fileId = "b16ec6b2-f807-46c9-9097-3be3590d0583"
$mdName = "$fileId.md"

# ---------------------------------------------------------------------
# Overview sheet -> new row 7
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(7, 1).Value = $mdName
$ov.Cells.Item(7, 2).Value = "Ready for handoff"
$ov.Cells.Item(7, 3).Value = "Ready for handoff"
$ov.Cells.Item(7, 4).Value = "2016-44-17 02:44:10"

$ov.Hyperlinks.Add(
    $ov.Cells.Item(7, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/512716f1922cca79adf8bf5c772ca61e9d66fa82/e2e/$mdName",
    "",
    "",
    $mdName
)

# ---------------------------------------------------------------------
# zh-cn sheet -> new row 7
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$fileId.7292d455f07c509ac44c124e633be620a55a2b77.zh-cn.xlf"

$zh.Cells.Item(7, 1).Value = $mdName
$zh.Cells.Item(7, 2).Value = ".md"
$zh.Cells.Item(7, 3).Value = "Ready for handoff"
$zh.Cells.Item(7, 4).Value = $zhXlf
$zh.Cells.Item(7, 5).Value = "2016-03-17 02:44:02"
$zh.Cells.Item(7, 8).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(7, 9).Value = "Include"

$zh.Hyperlinks.Add(
    $zh.Cells.Item(7, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/512716f1922cca79adf8bf5c772ca61e9d66fa82/e2e/$mdName",
    "",
    "",
    $mdName
)
$zh.Hyperlinks.Add(
    $zh.Cells.Item(7, 2),
    "https://github.com/OpenLocalizationTest/oltest/blob/512716f1922cca79adf8bf5c772ca61e9d66fa82/e2e/$mdName",
    "",
    "",
    ".md"
)
$zh.Hyperlinks.Add(
    $zh.Cells.Item(7, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/00000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
)

# ---------------------------------------------------------------------
# de-de sheet -> new row 7
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deXlf = "$fileId.7292d455f07c509ac44c124e633be620a55a2b77.de-de.xlf"

$de.Cells.Item(7, 1).Value = $mdName
$de.Cells.Item(7, 2).Value = ".md"
$de.Cells.Item(7, 3).Value = "Ready for handoff"
$de.Cells.Item(7, 4).Value = $deXlf
$de.Cells.Item(7, 5).Value = "2016-03-17 02:44:10"
$de.Cells.Item(7, 8).Value = "0001-01-01 00:00:00"
$de.Cells.Item(7, 9).Value = "Include"

$de.Hyperlinks.Add(
    $de.Cells.Item(7, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/512716f1922cca79adf8bf5c772ca61e9d66fa82/e2e/$mdName",
    "",
    "",
    $mdName
)
$de.Hyperlinks.Add(
    $de.Cells.Item(7, 2),
    "https://github.com/OpenLocalizationTest/oltest/blob/512716f1922cca79adf8bf5c772ca61e9d66fa82/e2e/$mdName",
    "",
    "",
    ".md"
)
$de.Hyperlinks.Add(
    $de.Cells.Item(7, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/00000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
)

Write-Host "Added handoff row for $fileId to Overview, zh-cn, de-de sheets"
